$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 previously held the text "R40" (shared string index 15).
# The edit replaces its content with the text "1" (a new shared string,
# stored as text rather than a number - note the leading apostrophe,
# Excel's standard "force text" quote-prefix convention).
$ws.Range("B11").Formula = "'1"
